$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$metaWs = $wb.Worksheets.Item("Metadata")

# Version: 2.2.0-ballot -> 2.1.0
$metaWs.Range("B3").Value = "2.1.0"

# Date: 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
$metaWs.Range("B8").Value = "2025-12-19T08:44:55+00:00"

# --- Include ValueSet #0 sheet ---
$ws0 = $wb.Worksheets.Item("Include ValueSet #0")
# ValueSet URL: drop the "|2.1.0" version suffix
$ws0.Range("A2").Value = "https://hl7.fr/ig/fhir/core/ValueSet/fr-core-vs-encounter-type"

# --- Include ValueSet #2 sheet ---
$ws2 = $wb.Worksheets.Item("Include ValueSet #2")
# ValueSet URL: drop the "|20250624152100" version suffix
$ws2.Range("A2").Value = "https://smt.esante.gouv.fr/fhir/ValueSet/jdv-type-evenement-ssiad-cisis"
